# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (referenced by the notes master)
#   ppt/theme/theme2.xml -> "Integral" colours     (referenced by the slide master
#                                                     and the presentation's main theme)
# The authored commit swaps the content of the two theme parts: the deck's
# visible design (slide master / presentation theme, i.e. theme2.xml) becomes
# the plain "Office Theme" colour scheme, while theme1.xml becomes "Integral".
#
# The slide master's theme colour scheme (the 12-slot a:clrScheme used by
# theme2.xml) is reachable and editable through the PowerPoint object model
# via Slide.ThemeColorScheme, so drive the swap through that.

function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target palette: the stock "Office Theme" colour scheme (12 slots, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order).
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToComRgb($officeThemeColors[$i - 1])
}

# Best-effort: also try to rename the colour scheme / design to match the
# stock "Office Theme" naming (no-op on hosts that don't expose a setter
# for this).
try { $tcs.Name = "Office" } catch {}
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
